# Daily attendance processing - 2025-11-29 10:51:02
#
# Normalizes the "Recorded By" column (G) on the "Session Analysis Results"
# sheet: for a fixed set of known author-string combinations, the order of
# the comma-separated recorder names is rotated so the last name listed is
# moved to the front (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact old-value -> new-value mapping observed for the "Recorded By" column.
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G = "Recorded By" (column index 7); data starts on row 2 (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
